$wb = $excel.ActiveWorkbook
$staff = $wb.Worksheets.Item("Staff")
$scenario = $wb.Worksheets.Item("Scenario")

# 1. The leave-factor helper column (F) used to pull its value from A35 via a
#    formula. The helper block at the bottom of the sheet is going away, so
#    bake the computed constant straight into F2:F10 instead of referencing it.
$staff.Range("F2:F10").Value = 0.80769199999999997

# 2. Drop the now-unused SUM helper cell and blank spacer row.
$staff.Range("L17").ClearContents()

# 3. Remove the old "LIMT oncall loss" total and the Leave-weeks / Weeks-per-year
#    / Leave-factor / On-call-factor explanation block that lived at the
#    bottom of the sheet (rows 18, 30, 32-36 in the original layout). Delete
#    from the bottom up so earlier row numbers stay valid.
$staff.Rows("36:36").Delete()
$staff.Rows("35:35").Delete()
$staff.Rows("34:34").Delete()
$staff.Rows("33:33").Delete()
$staff.Rows("32:32").Delete()
$staff.Rows("30:30").Delete()
$staff.Rows("18:18").Delete()

$staff.Range("D28").Select()

# 4. Tidy up the Scenario sheet row heights now that the descriptive text in
#    column A is shorter (no more wrapped 45pt-tall rows).
$scenario.Rows("3:3").RowHeight = 15
$scenario.Rows("4:4").RowHeight = 30
$scenario.Rows("5:5").RowHeight = 15
